# Update Work Week and Social Spending
# Updates GDP per Capita data values in the "Data" sheet and appends
# six new yearly rows (2011-2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Data")

# --- 1. Update existing data points (column E) with revised GDP figures ---
# Each entry is (row, newTextValue). Values are written as TEXT (matching
# the original workbook's storage of column E as shared strings), and
# ClearFormats() strips the transient quote-prefix formatting Excel applies
# when a numeric-looking string is entered, so the cell keeps the sheet's
# default (unstyled) appearance.
$updates = @(
    @(2, "931"),
    @(52, "995"),
    @(84, "1073"),
    @(85, "1282"),
    @(86, "1082"),
    @(87, "1130"),
    @(88, "1151"),
    @(89, "1191"),
    @(90, "1235"),
    @(91, "1235"),
    @(92, "1393"),
    @(93, "1455"),
    @(94, "1452"),
    @(95, "1575"),
    @(96, "1517"),
    @(97, "1395"),
    @(98, "1599"),
    @(99, "1830"),
    @(100, "2050"),
    @(101, "1932"),
    @(102, "2055"),
    @(103, "1969"),
    @(104, "2163"),
    @(105, "1992"),
    @(106, "2133"),
    @(107, "2109"),
    @(108, "2174"),
    @(109, "2177"),
    @(110, "2219"),
    @(111, "2252"),
    @(112, "2203"),
    @(113, "2168"),
    @(114, "2212"),
    @(115, "2181"),
    @(116, "2149"),
    @(117, "1956"),
    @(118, "2174"),
    @(119, "2276"),
    @(120, "2295"),
    @(121, "2404"),
    @(122, "2402"),
    @(128, "1030"),
    @(129, "1395"),
    @(130, "1583"),
    @(131, "1634"),
    @(132, "1706"),
    @(133, "1835"),
    @(134, "1890"),
    @(135, "1999"),
    @(136, "2085"),
    @(137, "2165"),
    @(138, "2248"),
    @(139, "2299"),
    @(140, "2308"),
    @(141, "2393"),
    @(142, "2353"),
    @(143, "2410"),
    @(144, "2450"),
    @(145, "2542"),
    @(146, "2550"),
    @(147, "2603"),
    @(148, "2636"),
    @(149, "2694"),
    @(150, "2745"),
    @(151, "2789"),
    @(152, "2812"),
    @(153, "2882"),
    @(154, "2954"),
    @(155, "3131"),
    @(156, "3154"),
    @(157, "3241"),
    @(158, "3430"),
    @(159, "3524"),
    @(160, "3606"),
    @(161, "3703"),
    @(162, "3787"),
    @(163, "3819"),
    @(164, "3859"),
    @(165, "3837"),
    @(166, "3469"),
    @(167, "3135"),
    @(168, "3161"),
    @(169, "3218"),
    @(170, "3355"),
    @(171, "3481"),
    @(172, "3502"),
    @(173, "3432.58275087477"),
    @(174, "3399.24261795966"),
    @(175, "3427.12943028549"),
    @(176, "3529.00733468657"),
    @(177, "3637.62556131435"),
    @(178, "3790.77878882303"),
    @(179, "3928.01438143073"),
    @(180, "3852.63274317198"),
    @(181, "3917.08864311008"),
    @(182, "4033.56089188223"),
    @(183, "4090.21832899351"),
    @(184, "4181.37574820598"),
    @(185, "4333.21381169289"),
    @(186, "4562.50267556432"),
    @(187, "4721.00072707225"),
    @(188, "4908.6168623321"),
    @(189, "5172.66187101868"),
    @(190, "5328.68244964514"),
    @(191, "5338.30624330191"),
    @(192, "5694.04672885807")
)

foreach ($u in $updates) {
    $row = $u[0]
    $val = $u[1]
    $cell = $ws.Cells.Item($row, 5)
    $cell.Formula = "'" + $val
    $cell.ClearFormats()
}

# --- 2. Append six new rows for years 2011-2016 ---
$newRows = @(
    @(193, 2011, "5851"),
    @(194, 2012, "6144"),
    @(195, 2013, "6472"),
    @(196, 2014, "6763"),
    @(197, 2015, "7047"),
    @(198, 2016, "7410")
)

foreach ($nr in $newRows) {
    $row = $nr[0]
    $year = $nr[1]
    $val = $nr[2]

    $ws.Cells.Item($row, 1).Value = 608
    $ws.Cells.Item($row, 2).Value = "Philippines"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year

    $cell = $ws.Cells.Item($row, 5)
    $cell.Formula = "'" + $val
    $cell.ClearFormats()
}
